# Generate Report for handoff
# - Updates the Status string shared across sheets ("Handoff transform failed" -> "Ready for handoff")
# - Records the first successful handoff (file + datetime) for zh-cn and de-de
# - Flips the Handoff Reason from "Ignored" to "Include" for the handed-off file

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$commitHash = "de790bef807620a1b01868421b7609cd4fe9e29d"
$baseName = "796ffe9a-282e-49df-a2e5-aa11cf21d28a"
$repoBlob = "https://github.com/OpenLocalizationTest/oltest/blob/e80369f8ee1f9f8920bdd657e7c75271f99d858a"

function Set-HandoffRow($SheetName, $HandoffFileName, $HandoffDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status: "Handoff transform failed" -> "Ready for handoff"
    $ws.Range("B2").Value = $newStatus

    # Latest Handoff File (column C) - newly populated with a hyperlink
    $ws.Range("C2").Value = $HandoffFileName
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBlob/loc/$HandoffFileName", "", "", $HandoffFileName) | Out-Null

    # Latest Handoff Datetime (column D)
    $ws.Range("D2").Value = $HandoffDateTime

    # Handoff Reason (column H): "Ignored" -> "Include"
    $ws.Range("H2").Value = "Include"
}

# zh-cn sheet
Set-HandoffRow "zh-cn" "$baseName.$commitHash.zh-cn.xlf" "2016-01-18 04:08:25"

# de-de sheet
Set-HandoffRow "de-de" "$baseName.$commitHash.de-de.xlf" "2016-01-18 04:08:37"

# Overview sheet mirrors the Status column for each locale
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
